# Apply the cryptocurrency price/volume refresh described in the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.069.48"
$ws.Range("E2").Value = "  +0.08%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.833.49"
$ws.Range("E3").Value = "  +0.17%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9978"
$ws.Range("E4").Value = "  -0.08%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.05"
$ws.Range("E5").Value = "  +1.58%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6344"
$ws.Range("E6").Value = "  +1.20%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9997"
$ws.Range("E7").Value = "  -0.01%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07559"
$ws.Range("E8").Value = "  -0.45%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2950"
$ws.Range("E9").Value = "  +1.26%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.23"
$ws.Range("E10").Value = "  +1.96%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07717"
$ws.Range("E11").Value = "  +1.05%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.835.14"

# Row 13
$ws.Range("E13").Value = "  +1.07%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6737"
$ws.Range("E14").Value = "  +1.26%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "83.27"
$ws.Range("E15").Value = "  +1.22%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000009590"
$ws.Range("E16").Value = "  +5.04%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.103"
$ws.Range("E17").Value = "  +2.01%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "29.103.31"
$ws.Range("E18").Value = "  +0.41%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.63"
$ws.Range("E19").Value = "  +2.40%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "227.49"
$ws.Range("E20").Value = "  +1.31%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9984"
$ws.Range("E21").Value = "  -0.15%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.182"
$ws.Range("E22").Value = "  -0.21%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.9992"
$ws.Range("E23").Value = "  -0.09%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "160.10"
$ws.Range("E24").Value = "  +0.18%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1431"
$ws.Range("E25").Value = "  +5.06%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.551"
$ws.Range("E26").Value = "  +1.63%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.505"
$ws.Range("E28").Value = "  +0.67%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.157"
$ws.Range("E29").Value = "  +2.70%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.079"
$ws.Range("E30").Value = "  +1.18%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.05473"
$ws.Range("E31").Value = "  +5.34%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.202"
$ws.Range("E32").Value = "  -0.23%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.861"
$ws.Range("E33").Value = "  +0.81%  "

# Row 34
$ws.Range("E34").Value = "  +2.18%  "

# Row 35
$ws.Range("E35").Value = "  -0.98%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.659"
$ws.Range("E36").Value = "  +1.67%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.248.09"
$ws.Range("E37").Value = "  -2.47%  "

# Row 38
$ws.Range("E38").Value = "  -0.06%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01790"
$ws.Range("E39").Value = "  +0.03%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.640"
$ws.Range("E40").Value = "  +2.08%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9067"
$ws.Range("E41").Value = "  +1.65%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9997"
$ws.Range("E42").Value = "  +0.00%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "101.77"
$ws.Range("E43").Value = "  +0.20%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.985.18"
$ws.Range("E44").Value = "  +0.41%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00000000123"
$ws.Range("E45").Value = "  +3.13%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "65.25"
$ws.Range("E46").Value = "  +2.14%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5101"
$ws.Range("E47").Value = "  -0.14%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4084"
$ws.Range("E48").Value = "  +2.58%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.056"
$ws.Range("E49").Value = "  +2.57%  "

# Row 50
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05793"
$ws.Range("E50").Value = "  +0.95%  "

# Row 51
$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.652"
$ws.Range("E51").Value = "  +0.27%  "
